$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds formatted numeric-looking text (e.g. thousand-separator
# dotted strings like "58.287.14") in the source data. Force the column to Text
# first so Excel does not reinterpret values such as "559.26" as a float and
# round-trip them with binary floating point noise.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "58.287.14"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "2.995.04"
$ws.Range("E3").Value = "  +3.47%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "559.26"
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("D6").Value = "135.43"
$ws.Range("E6").Value = "  +11.20%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "0.517"
$ws.Range("E8").Value = "  +4.52%  "
$ws.Range("D9").Value = "2.991.33"
$ws.Range("E9").Value = "  +3.66%  "
$ws.Range("E10").Value = "  +4.45%  "
$ws.Range("D11").Value = "4.88"
$ws.Range("E11").Value = "  +2.75%  "
$ws.Range("D12").Value = "0.454"
$ws.Range("E12").Value = "  +4.97%  "
$ws.Range("D13").Value = "0.0000224"
$ws.Range("E13").Value = "  +6.17%  "
$ws.Range("D14").Value = "33.26"
$ws.Range("E14").Value = "  +5.31%  "
$ws.Range("D15").Value = "0.123"
$ws.Range("E15").Value = "  +3.27%  "
$ws.Range("D16").Value = "3.486.45"
$ws.Range("E16").Value = "  +3.24%  "
$ws.Range("D17").Value = "7.06"
$ws.Range("E17").Value = "  +9.51%  "
$ws.Range("D18").Value = "2.979.52"
$ws.Range("E18").Value = "  +2.97%  "
$ws.Range("D19").Value = "58.140.04"
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("D20").Value = "422.38"
$ws.Range("E20").Value = "  +3.41%  "
$ws.Range("D21").Value = "13.60"
$ws.Range("E21").Value = "  +6.42%  "
$ws.Range("D22").Value = "0.713"
$ws.Range("E22").Value = "  +9.44%  "
$ws.Range("D23").Value = "13.49"
$ws.Range("E23").Value = "  +7.12%  "
$ws.Range("D24").Value = "7.10"
$ws.Range("E24").Value = "  +5.55%  "
$ws.Range("D25").Value = "80.47"
$ws.Range("E25").Value = "  +4.76%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "0.997"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").Value = "2.52"
$ws.Range("E28").Value = "  +2.40%  "
$ws.Range("D29").Value = "2.06"
$ws.Range("E29").Value = "  +8.38%  "
$ws.Range("D30").Value = "7.61"
$ws.Range("E30").Value = "  +6.23%  "
$ws.Range("D31").Value = "25.81"
$ws.Range("E31").Value = "  +4.98%  "
$ws.Range("D32").Value = "6.03"
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("D33").Value = "0.0978"
$ws.Range("E33").Value = "  +2.99%  "
$ws.Range("D34").Value = "5.76"
$ws.Range("E34").Value = "  +7.83%  "
$ws.Range("D35").Value = "0.962"
$ws.Range("E35").Value = "  +7.08%  "
$ws.Range("D36").Value = "2.10"
$ws.Range("E36").Value = "  +3.72%  "
$ws.Range("D37").Value = "0.0₃0720"
$ws.Range("E37").Value = "  +16.90%  "
$ws.Range("D38").Value = "8.91"
$ws.Range("E38").Value = "  +6.79%  "
$ws.Range("D39").Value = "48.48"
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("D40").Value = "2.74"
$ws.Range("E40").Value = "  +17.29%  "
$ws.Range("D41").Value = "393.69"
$ws.Range("E41").Value = "  +9.50%  "
$ws.Range("D42").Value = "0.0350"
$ws.Range("E42").Value = "  +1.90%  "
$ws.Range("D43").Value = "0.108"
$ws.Range("E43").Value = "  +2.44%  "
$ws.Range("D44").Value = "2.731.21"
$ws.Range("E44").Value = "  +4.67%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "125.33"
$ws.Range("E46").Value = "  +6.51%  "
$ws.Range("D47").Value = "0.242"
$ws.Range("E47").Value = "  +6.36%  "
$ws.Range("D48").Value = "2.01"
$ws.Range("E48").Value = "  +4.31%  "
$ws.Range("D49").Value = "0.109"
$ws.Range("E49").Value = "  +2.83%  "
$ws.Range("D50").Value = "23.22"
$ws.Range("E50").Value = "  +4.70%  "
$ws.Range("D51").Value = "2.02"
$ws.Range("E51").Value = "  +4.27%  "
